$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card9")

# Fix the M1 header text: remove trailing space from "Event "
$ws.Range("M1").Value = "Event"

# Add the new header "Correction " (with trailing space) in N1, copying the
# header style (bold font, border, centered/top alignment) from M1
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("N1").Value = "Correction "

# Data rows: M2:M12 previously empty -> now hold the text "nan";
# N2:N12 are new, left empty (but present as real cells on the sheet).
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"

    # Touch the new N cell (without giving it a value) so it materialises
    # as a real, empty cell on the sheet rather than being left absent.
    $ws.Cells.Item($r, 14).Font.Bold = $false
}

$excel.CutCopyMode = 0
